$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metadata sheet: bump the generation Date property.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-01-23T08:28:04+00:00"

# ---------------------------------------------------------------------------
# 2) Mapping Table 0 sheet: fix the "consumable" separator typo
#    (dot -> colon) on the LM -> CDA "vaccin" mapping row.
# ---------------------------------------------------------------------------
$table0 = $wb.Worksheets.Item("Mapping Table 0")
$table0.Range("D12").Value = "FRCDAVaccination.consumable:FRCDAProduitDeSante"

# ---------------------------------------------------------------------------
# 3) Mapping Table 1 sheet: insert a new "doseQuantity" equivalence row
#    right after the routeCode/route row (was row 11, now everything from
#    old row 11 onward shifts down by one).
# ---------------------------------------------------------------------------
$table1 = $wb.Worksheets.Item("Mapping Table 1")

# Grow the table by one row, copying the formatting of the current last
# data row (16) down onto the new last row (17) so no new cell style gets
# minted.
$table1.Range("A16:E16").Copy()
$table1.Range("A17:E17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Shift the existing mapping rows 11..16 down into 12..17, walking from the
# bottom up so a row is never overwritten before it has been copied.
for ($r = 16; $r -ge 11; $r--) {
    $dest = $r + 1
    $table1.Range("A$r").Copy()
    $table1.Range("A$dest").PasteSpecial(-4163)
    $table1.Range("C$r").Copy()
    $table1.Range("C$dest").PasteSpecial(-4163)
    $table1.Range("D$r").Copy()
    $table1.Range("D$dest").PasteSpecial(-4163)
}
$excel.CutCopyMode = $false

# Write the new doseQuantity equivalence mapping into the freed-up row 11.
$table1.Range("A11").Value = "FRCDAVaccination.doseQuantity"
$table1.Range("C11").Value = "equivalent"
$table1.Range("D11").Value = "FRImmunizationDocument.doseQuantity"

Write-Output "edit applied"
